$d = $word.ActiveDocument

# 1. Replace the title text
$d.Content.Find.Execute("2.2 - Debate I", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Placeholder - Check Back Later", 2)

# 2. Remove the trailing " " and ":::" runs from the last paragraph of the
#    "Additional Resources" bullet (table cell). Find and delete the
#    " :::" text that follows "...general edification later."
$d.Content.Find.Execute(" :::", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)
